# Portfolio workbook update:
#  - watchlist: refresh prices/pct_change/rsi and re-sort rows by rsi ascending
#  - stocks: manual sell of all WMT shares -> remove its row, refresh current
#    price/value/performance/rsi for the remaining holdings
#  - portfolio: refreshed CASH/STOCKS/TOTAL after the sale
#  - trades: append the new WMT sell trade
#  - summary: append the corresponding portfolio snapshot row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. watchlist sheet - rewrite rows 2..31 (ticker, price, pct_change, rsi)
# ---------------------------------------------------------------------------
$wl = $wb.Worksheets.Item("watchlist")

$watchlistRows = @(
    @("CSCO", 42.09, -0.9647, 23.003194888179),
    @("MSFT", 210.28, 0.6606, 38.76525484565688),
    @("PFE", 38.35, 0.762, 50),
    @("WBA", 41.17, -1.7657, 50.7718696397942),
    @("TRV", 115.91, -2.2681, 54.83045425463848),
    @("JNJ", 148.99, 0.5059, 58.0550098231828),
    @("GS", 203.07, -2.3561, 58.23634735899732),
    @("INTC", 48.93, 0.0818, 58.2474226804124),
    @("IBM", 124.44, -0.6626, 60.45673076923076),
    @("XOM", 42.64, -1.2963, 61.39359698681731),
    @("BA", 172.01, -3.4086, 61.82038834951454),
    @("AXP", 97.55, -2.8483, 62.43025418474887),
    @("JPM", 99.70999999999999, -2.6365, 63.87394312067643),
    @("AAPL", 458.43, -0.2611, 65.81288757077375),
    @("MRK", 84.76000000000001, 1.5333, 69.82142857142848),
    @("CAT", 138.72, -0.886, 71.26225490196079),
    @("PG", 135.5, 0.2961, 72.22808870116161),
    @("VZ", 58.78, -0.017, 72.51461988304095),
    @("WMT", 135.6, 2.2624, 73.18397827562788),
    @("V", 199.43, 1.4188, 75.87392550143264),
    @("UNH", 320.51, -0.9855, 76.50564617314927),
    @("NKE", 105.66, -0.7235, 79.06423473433789),
    @("DIS", 129.37, -0.8887, 79.69890510948909),
    @("CVX", 90.77, 0.4649, 81.97424892703857),
    @("RTX", 62.77, -1.3671, 82.42009132420087),
    @("DOW", 44.46, -0.9137999999999999, 85.07223113964695),
    @("KO", 48.21, -0.4954, 86.0377358490567),
    @("MMM", 164.71, -0.8368, 91.86182669789234),
    @("HD", 288.24, 2.741, 94.01436552274541),
    @("MCD", 208.67, 0.7922, 96.2410887880751)
)

$r = 2
foreach ($row in $watchlistRows) {
    $wl.Cells($r, 1).Value = $row[0]
    $wl.Cells($r, 2).Value = $row[1]
    $wl.Cells($r, 3).Value = $row[2]
    $wl.Cells($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. stocks sheet - sell all WMT shares: drop its row, refresh the rest
# ---------------------------------------------------------------------------
$st = $wb.Worksheets.Item("stocks")

# find & delete the WMT row
$wmtRow = 0
for ($i = 2; $i -le 8; $i = $i + 1) {
    if ($st.Cells($i, 1).Value2 -eq "WMT") {
        $wmtRow = $i
    }
}
$st.Rows($wmtRow).Delete()

$stocksRows = @(
    @("XOM", 42.64, 511.68, -7.264, 61.39359698681731),
    @("INTC", 48.93, 880.74, -3.2813, 58.2474226804124),
    @("JNJ", 148.99, 595.96, 1.4642, 58.0550098231828),
    @("TRV", 115.91, 463.64, 2.4845, 54.83045425463848),
    @("PFE", 38.35, 536.9, 1.4819, 50),
    @("CSCO", 42.09, 589.26, -0.9647, 23.003194888179)
)

$r = 2
foreach ($row in $stocksRows) {
    # sanity-check row alignment by ticker before overwriting
    if ($st.Cells($r, 1).Value2 -eq $row[0]) {
        $st.Cells($r, 3).Value = $row[1]
        $st.Cells($r, 5).Value = $row[2]
        $st.Cells($r, 6).Value = $row[3]
        $st.Cells($r, 7).Value = $row[4]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. portfolio sheet - refreshed CASH / STOCKS / TOTAL
# ---------------------------------------------------------------------------
$pf = $wb.Worksheets.Item("portfolio")
$pf.Cells(2, 2).Value = 7431.629999999999
$pf.Cells(3, 2).Value = 3578.18
$pf.Cells(4, 2).Value = 11009.81

# ---------------------------------------------------------------------------
# 4. trades sheet - append the manual WMT sell
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("trades")
$trLastRow = $tr.Cells($tr.Rows.Count, 2).End(-4162).Row
$trNewRow = $trLastRow + 1

$tr.Cells($trLastRow, 1).Copy($tr.Cells($trNewRow, 1))
$tr.Cells($trNewRow, 1).Value = 67
$tr.Cells($trNewRow, 2).Value = "17/08/2020 17:56:45"
$tr.Cells($trNewRow, 3).Value = "WMT"
$tr.Cells($trNewRow, 4).Value = "sell"
$tr.Cells($trNewRow, 5).Value = 4
$tr.Cells($trNewRow, 6).Value = 542.4

# ---------------------------------------------------------------------------
# 5. summary sheet - append the matching portfolio snapshot
# ---------------------------------------------------------------------------
$sm = $wb.Worksheets.Item("summary")
$smLastRow = $sm.Cells($sm.Rows.Count, 1).End(-4162).Row
$smNewRow = $smLastRow + 1

$sm.Cells($smLastRow, 1).Copy($sm.Cells($smNewRow, 1))
$sm.Cells($smNewRow, 1).Value = "17/08/2020 17:56:45"
$sm.Cells($smNewRow, 2).Value = 7431.629999999999
$sm.Cells($smNewRow, 3).Value = 3578.18
$sm.Cells($smNewRow, 4).Value = 11009.81
